# Update existing TPM-derived values (rows 2-4) and append a new row 5
# for the "Resolving-Mac" target cluster, per the refreshed TPM scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value = 0.06624833333333334
$ws.Range("H2").Value = 0.198745
$ws.Range("M2").Value = 0.138766
$ws.Range("N2").Value = 0.416298
$ws.Range("O2").Value = 0.01356925767068476
$ws.Range("P2").Value = 0.01356925767068476
$ws.Range("Q2").Value = 0.009193016223333334
$ws.Range("R2").Value = 0.08273714601
$ws.Range("S2").Value = 0.01356925767068476
$ws.Range("T2").Value = 0.01356925767068476

# --- Row 3 updates ---
$ws.Range("G3").Value = 0.06624833333333334
$ws.Range("H3").Value = 0.198745
$ws.Range("O3").Value = 0.00529617548814441
$ws.Range("P3").Value = 0.005296175488144411
$ws.Range("Q3").Value = 0.003588098064444445
$ws.Range("R3").Value = 0.03229288258
$ws.Range("S3").Value = 0.00529617548814441
$ws.Range("T3").Value = 0.005296175488144411

# --- Row 4 updates ---
$ws.Range("G4").Value = 0.06624833333333334
$ws.Range("H4").Value = 0.198745
$ws.Range("M4").Value = 10.002366
$ws.Range("N4").Value = 30.007098
$ws.Range("O4").Value = 0.9780831152479456
$ws.Range("P4").Value = 0.9780831152479456
$ws.Range("Q4").Value = 0.66264007689
$ws.Range("R4").Value = 5.96376069201
$ws.Range("S4").Value = 0.9780831152479456
$ws.Range("T4").Value = 0.9780831152479456

# --- New row 5: MuSCs -> Cntn1/Ptprz1 -> Resolving-Mac ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Cntn1"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.06624833333333334
$ws.Range("H5").Value = 0.198745
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03120566666666667
$ws.Range("N5").Value = 0.09361700000000001
$ws.Range("O5").Value = 0.003051451593225274
$ws.Range("P5").Value = 0.003051451593225274
$ws.Range("Q5").Value = 0.002067323407222223
$ws.Range("R5").Value = 0.018605910665
$ws.Range("S5").Value = 0.003051451593225274
$ws.Range("T5").Value = 0.003051451593225274
